# Auto-generated edit script: updates crypto price/name/link/volume cells
# per the commit diff (GitHub Actions symbol-list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''269.44'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').Value = '''22.82'
$ws.Range('D3').ClearFormats()
$ws.Range('D4').Value = '''6.333'
$ws.Range('D4').ClearFormats()
$ws.Range('D5').Value = '''0.06188'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').Value = '''3.641'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').Value = '''6.678'
$ws.Range('D7').ClearFormats()
$ws.Range('D8').Value = '''1.389'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').Value = '''0.8301'
$ws.Range('D9').ClearFormats()
$ws.Range('D11').Value = '''0.1605'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').Value = '''0.08229'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').Value = '''0.03465'
$ws.Range('D13').ClearFormats()
$ws.Range('D14').Value = '''0.03228'
$ws.Range('D14').ClearFormats()
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').Value = '''0.09318'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').Value = '''3.849'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').Value = '''0.001641'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').Value = '''0.04738'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D19').Value = '''0.006334'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D20').Value = '''0.005669'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D21').Value = '''0.001077'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D22').Value = '''0.0001499'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'LEO'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D23').Value = '''3.720'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D24').Value = '''2.413'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('B25').Value = 'BitpandaEcosystemToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D25').Value = '''0.3340'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'
$ws.Range('B26').Value = 'ProBitToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D26').Value = '''0.1239'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '25ProBitTokenPROB'
$ws.Range('D27').Value = '''0.0002703'
$ws.Range('D27').ClearFormats()
$ws.Range('D41').Value = '''0.006954'
$ws.Range('D41').ClearFormats()
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = '''0.1162'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').Value = '''0.003347'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').Value = '''0.01154'
$ws.Range('D44').ClearFormats()
$ws.Range('D47').Value = '''0.00000000750'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').Value = '''0.9199'
$ws.Range('D48').ClearFormats()
$ws.Range('B49').Value = 'CryptobidCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'
$ws.Range('D49').Value = '''0.00001400'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '48CryptobidCoinCBCWorstin24h'
$ws.Range('B50').Value = 'BOLO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D50').Value = '''0.002270'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '49BOLOBOLO'
$ws.Range('D51').Value = '''0.01240'
$ws.Range('D51').ClearFormats()
